# Update "想去人数" (interested-count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 597
$ws1.Range("F4").Value = 33
$ws1.Range("F6").Value = 356
$ws1.Range("F7").Value = 1753

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 597
$ws4.Range("F4").Value = 33
$ws4.Range("F6").Value = 356
$ws4.Range("F11").Value = 1753
